$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Администратор, менеджер по подбору персонала, рекрутер, помощник руководителя"
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = "Томская область, г Томск"
$ws.Range("D2").Value = "Высшее-бакалавриат"
$ws.Range("E2").Value = 25000

# Row 3
$ws.Range("A3").Value = "учитель русского языка, учитель английского языка"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "г Москва"
$ws.Range("D3").Value = "Высшее-бакалавриат"
$ws.Range("E3").Value = 50000

# Row 4
$ws.Range("A4").Value = "инженер-программист"
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = "Московская область, г Дубна"
$ws.Range("D4").Value = "Высшее-бакалавриат"
$ws.Range("E4").Value = 150000

# Row 5
$ws.Range("A5").Value = "Педагог, репетитор, администратор, менеджер"
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = "Московская область, г Балашиха"
$ws.Range("D5").Value = "Среднее общее"
$ws.Range("E5").Value = 20000

# Row 6
$ws.Range("A6").Value = "менеджер, преподаватель"
$ws.Range("B6").Value = 15
$ws.Range("C6").Value = "Самарская область, г Самара"
$ws.Range("D6").Value = "Высшее-бакалавриат"
$ws.Range("E6").Value = 80000

# Row 7
$ws.Range("A7").Value = "Програмист, HTML-верстальщик, Wordpress-разработчик"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Краснодарский край, г Новороссийск, Натухаевская станица"
$ws.Range("D7").Value = "Среднее общее"
$ws.Range("E7").Value = 30000

# Row 8
$ws.Range("A8").Value = "Педагог-библиотекарь, менеджер научно-исследовательской деятельности, методист, заведующий библиотеки"
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = "Тюменская область, г Тобольск"
$ws.Range("D8").Value = "Высшее-бакалавриат"
$ws.Range("E8").Value = 25000

# Row 9
$ws.Range("A9").Value = "Программист-экономист, 1с-программист"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = "Владимирская область, г Владимир"
$ws.Range("D9").Value = "Высшее-бакалавриат"
$ws.Range("E9").Value = 40000

# Row 10
$ws.Range("A10").Value = "Педагог, завуч, администратор автоматизированных баз данных, библиограф, библиотекарь, психолог, арт-терапевт"
$ws.Range("B10").Value = 23
$ws.Range("C10").Value = "Московская область"
$ws.Range("D10").Value = "Высшее-бакалавриат"
$ws.Range("E10").Value = 40000

# Row 11
$ws.Range("A11").Value = "Педагог-организатор ОБЖ, ОТ и ПБ, ГО и ЧС"
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "г Москва"
$ws.Range("D11").Value = "Высшее-бакалавриат"
$ws.Range("E11").Value = 80000
